$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin name, link, price, and volume cells to refresh the crypto listing
$ws.Range("D2").Value = "36.627.04"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "2.047.80"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -4.79%  "
$ws.Range("E9").Value = "  +7.31%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("E12").Value = "  -3.53%  "
$ws.Range("E13").Value = "  +9.31%  "
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "2.338.97"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").Value = "2.044.90"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "36.551.32"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").Value = "0.0₃0857"
$ws.Range("E21").Value = "  -3.35%  "
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("E23").Value = "  -3.60%  "
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("E27").Value = "  -6.96%  "
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("E31").Value = "  -8.01%  "
$ws.Range("E32").Value = "  +4.49%  "
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("E34").Value = "  -6.92%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("E40").Value = "  +3.46%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E41").Value = "  -3.27%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E42").Value = "  -3.38%  "
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("E45").Value = "  -4.75%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.380.72"
$ws.Range("E47").Value = "  +5.93%  "
$ws.Range("E48").Value = "  +11.67%  "
$ws.Range("E49").Value = "  +2.93%  "
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("E51").Value = "  +1.42%  "

# Cells whose new text values resemble numbers need an explicit text format
# so Excel preserves the exact string (e.g. trailing zeros) instead of
# reinterpreting them as numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.663"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.56"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "63.36"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.368"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0750"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.955"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0600"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0869"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0215"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "93.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0904"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.76"
$ws.Range("D51").Style = "Normal"
